$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.003374998590558752
$ws.Range("C2").Value = 0.0236459363664146

$ws.Range("B4").Value = 0.2307393176428526
$ws.Range("C4").Value = 0.1409991438301644

$ws.Range("B6").Value = 0.01146414169361983
$ws.Range("C6").Value = 0.04895474487947783

$ws.Range("B7").Value = 0.03083300495780911
$ws.Range("C7").Value = 0.05682242393025253

$ws.Range("B8").Value = 0.01122025094545002
$ws.Range("C8").Value = 0.05579377016739856
